$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51, pushing existing rows 51-52 down to 52-53.
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new price record.
$ws.Range("A51").Value = 7
$ws.Range("B51").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C51").Value = "Ñuble"
$ws.Range("D51").Value = 44461
$ws.Range("D51").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E51").Value = 16
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100108
$ws.Range("H51").Value = "Tropicales y subtropicales"
$ws.Range("I51").Value = 100108002
$ws.Range("J51").Value = "Mango"
$ws.Range("K51").Value = "Sin especificar"
$ws.Range("L51").Value = "Primera"
$ws.Range("M51").Value = 60
$ws.Range("N51").Value = 8500
$ws.Range("O51").Value = 9000
$ws.Range("P51").Value = 8750
$ws.Range("Q51").Value = "$/bandeja 4 kilos"
$ws.Range("R51").Value = "Brasil"
$ws.Range("S51").Value = 2188
$ws.Range("T51").Value = 4
